{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Title: \"Week 9\" -> \"Week 8: Core Exam\"\nparagraphs.items[0].insertText(\"Week 8: Core Exam\", Word.InsertLocation.replace);\n\n// 2) Subtitle: \"Data Wrangling and Visualization\" -> \"Practicum 1 Due Week 9\"\nparagraphs.items[1].insertText(\"Practicum 1 Due Week 9\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Author (\"Prof. Jack Reilly\") and Date (\"F2025\") paragraphs are unchanged.\n\n// 3) Replace the old \"Nothing to see here (yet)\" paragraph with two new\n//    sections (\"Readings\" and \"Submission\"), each a Heading 2 followed by\n//    a FirstParagraph-styled body paragraph.\nconst readingsHeading = paragraphs.items[4];\nreadingsHeading.insertText(\"Readings\", Word.InsertLocation.replace);\nconst readingsBody = readingsHeading.insertParagraph(\n  \"No readings or assignment due this week.\",\n  Word.InsertLocation.after\n);\nreadingsHeading.style = \"Heading 2\";\n\nconst submissionHeading = readingsBody.insertParagraph(\"Submission\", Word.InsertLocation.after);\nsubmissionHeading.style = \"Heading 2\";\nconst submissionBody = submissionHeading.insertParagraph(\n  \"Practicum 1 is due next week.\",\n  Word.InsertLocation.after\n);\nsubmissionBody.style = \"First Paragraph\";\nawait context.sync();\n\n// 4) Wrap \"Readings\" + its body paragraph in a \"readings\" bookmark, and\n//    \"Submission\" + its body paragraph in a \"submission\" bookmark.\nconst readingsRange = readingsHeading\n  .getRange(Word.RangeLocation.start)\n  .expandTo(readingsBody.getRange(Word.RangeLocation.end));\nreadingsRange.insertBookmark(\"readings\");\n\nconst submissionRange = submissionHeading\n  .getRange(Word.RangeLocation.start)\n  .expandTo(submissionBody.getRange(Word.RangeLocation.end));\nsubmissionRange.insertBookmark(\"submission\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Title: \"Week 9\" -> \"Week 8: Core Exam\"\n$d.Paragraphs(1).Range.Text = \"Week 8: Core Exam\"\n\n# 2) Subtitle: \"Data Wrangling and Visualization\" -> \"Practicum 1 Due Week 9\"\n$d.Paragraphs(2).Range.Text = \"Practicum 1 Due Week 9\"\n\n# Author (Prof. Jack Reilly) and Date (F2025) paragraphs are unchanged.\n\n# 3) Replace the old \"Nothing to see here (yet)\" paragraph with two new\n#    sections (\"Readings\" and \"Submission\"), each a Heading 2 followed by\n#    a FirstParagraph-styled body paragraph, wrapped in bookmarks.\n$pLast = $d.Paragraphs(5)\n$r = $pLast.Range\n$r.Text = \"Readings\"\n$r.InsertParagraphAfter()\n$pLast.Style = \"Heading 2\"\n\n$pBody1 = $d.Paragraphs(6)\n$pBody1.Range.Text = \"No readings or assignment due this week.\"\n$pBody1.Range.InsertParagraphAfter()\n\n$pHead2 = $d.Paragraphs(7)\n$pHead2.Range.Text = \"Submission\"\n$pHead2.Style = \"Heading 2\"\n$pHead2.Range.InsertParagraphAfter()\n\n$pBody2 = $d.Paragraphs(8)\n$pBody2.Range.Text = \"Practicum 1 is due next week.\"\n$pBody2.Style = \"First Paragraph\"\n\n# 4) Wrap \"Readings\" + its body paragraph in a \"readings\" bookmark, and\n#    \"Submission\" + its body paragraph in a \"submission\" bookmark.\n$readingsRange = $d.Range($d.Paragraphs(5).Range.Start, $d.Paragraphs(6).Range.End)\n$d.Bookmarks.Add(\"readings\", $readingsRange)\n\n$submissionRange = $d.Range($d.Paragraphs(7).Range.Start, $d.Paragraphs(8).Range.End)\n$d.Bookmarks.Add(\"submission\", $submissionRange)\n"}
